$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
